$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 25,9
$arr[0,0] = "model_6_4_24"
$arr[0,1] = 0.3083830318006148
$arr[0,2] = -0.01372893691086063
$arr[0,3] = 0.2805454215328287
$arr[0,4] = 0.1885902720498855
$arr[0,5] = 0.765416145324707
$arr[0,6] = 0.829413115978241
$arr[0,7] = 1.004480004310608
$arr[0,8] = 0.9117977023124695
$arr[1,0] = "model_6_4_20"
$arr[1,1] = 0.309036262977647
$arr[1,2] = -0.01338247709329088
$arr[1,3] = 0.2841236939266608
$arr[1,4] = 0.1908162146176458
$arr[1,5] = 0.7646932601928711
$arr[1,6] = 0.8291295170783997
$arr[1,7] = 0.9994841814041138
$arr[1,8] = 0.9092963933944702
$arr[2,0] = "model_6_4_23"
$arr[2,1] = 0.3090754762485025
$arr[2,2] = -0.01299818608712489
$arr[2,3] = 0.2839530333078509
$arr[2,4] = 0.190864720065346
$arr[2,5] = 0.7646498084068298
$arr[2,6] = 0.8288151621818542
$arr[2,7] = 0.9997223615646362
$arr[2,8] = 0.909241795539856
$arr[3,0] = "model_6_4_22"
$arr[3,1] = 0.3091756268427329
$arr[3,2] = -0.01321716159524144
$arr[3,3] = 0.2847251862190177
$arr[3,4] = 0.19123155680993
$arr[3,5] = 0.7645389437675476
$arr[3,6] = 0.8289943337440491
$arr[3,7] = 0.9986443519592285
$arr[3,8] = 0.9088295698165894
$arr[4,0] = "model_6_4_21"
$arr[4,1] = 0.3092442638602823
$arr[4,2] = -0.01321844769953473
$arr[4,3] = 0.2851110089223748
$arr[4,4] = 0.1914568583348523
$arr[4,5] = 0.7644630074501038
$arr[4,6] = 0.8289953470230103
$arr[4,7] = 0.9981056451797485
$arr[4,8] = 0.908576488494873
$arr[5,0] = "model_6_4_18"
$arr[5,1] = 0.3095232211170844
$arr[5,2] = -0.01328021366845755
$arr[5,3] = 0.2868826242920632
$arr[5,4] = 0.1924686001133882
$arr[5,5] = 0.764154314994812
$arr[5,6] = 0.8290459513664246
$arr[5,7] = 0.9956321716308594
$arr[5,8] = 0.9074395298957825
$arr[6,0] = "model_6_4_19"
$arr[6,1] = 0.3097514130877069
$arr[6,2] = -0.01383070360063687
$arr[6,3] = 0.2884219994874719
$arr[6,4] = 0.1931564449383407
$arr[6,5] = 0.7639017701148987
$arr[6,6] = 0.8294963240623474
$arr[6,7] = 0.9934829473495483
$arr[6,8] = 0.9066665768623352
$arr[7,0] = "model_6_4_17"
$arr[7,1] = 0.3107915320417168
$arr[7,2] = -0.01276072809831441
$arr[7,3] = 0.2936606414245554
$arr[7,4] = 0.1966317062240425
$arr[7,5] = 0.7627506852149963
$arr[7,6] = 0.8286207914352417
$arr[7,7] = 0.9861689805984497
$arr[7,8] = 0.9027613997459412
$arr[8,0] = "model_6_4_7"
$arr[8,1] = 0.3109292126275487
$arr[8,2] = -0.01291616180967892
$arr[8,3] = 0.2952249108119063
$arr[8,4] = 0.1974865815986456
$arr[8,5] = 0.7625983357429504
$arr[8,6] = 0.8287479877471924
$arr[8,7] = 0.9839849472045898
$arr[8,8] = 0.9018006324768066
$arr[9,0] = "model_6_4_15"
$arr[9,1] = 0.3109735114714482
$arr[9,2] = -0.0124882961751045
$arr[9,3] = 0.2945697432844298
$arr[9,4] = 0.1972683346548059
$arr[9,5] = 0.7625492811203003
$arr[9,6] = 0.8283979892730713
$arr[9,7] = 0.9848996996879578
$arr[9,8] = 0.9020459651947021
$arr[10,0] = "model_6_4_11"
$arr[10,1] = 0.3110550780031275
$arr[10,2] = -0.01207419810122778
$arr[10,3] = 0.2948625446699953
$arr[10,4] = 0.1975993880081425
$arr[10,5] = 0.7624589800834656
$arr[10,6] = 0.8280590772628784
$arr[10,7] = 0.9844908714294434
$arr[10,8] = 0.9016739130020142
$arr[11,0] = "model_6_4_12"
$arr[11,1] = 0.3110550780031275
$arr[11,2] = -0.01207419810122778
$arr[11,3] = 0.2948625446699953
$arr[11,4] = 0.1975993880081425
$arr[11,5] = 0.7624589800834656
$arr[11,6] = 0.8280590772628784
$arr[11,7] = 0.9844908714294434
$arr[11,8] = 0.9016739130020142
$arr[12,0] = "model_6_4_6"
$arr[12,1] = 0.3112762669840543
$arr[12,2] = -0.01284642248423262
$arr[12,3] = 0.2972002844038775
$arr[12,4] = 0.1986684103882895
$arr[12,5] = 0.7622142434120178
$arr[12,6] = 0.8286910057067871
$arr[12,7] = 0.9812270402908325
$arr[12,8] = 0.9004727005958557
$arr[13,0] = "model_6_4_5"
$arr[13,1] = 0.31132893390105
$arr[13,2] = -0.01270738536156046
$arr[13,3] = 0.2974703697017154
$arr[13,4] = 0.1988795746089219
$arr[13,5] = 0.7621559500694275
$arr[13,6] = 0.8285772204399109
$arr[13,7] = 0.9808499813079834
$arr[13,8] = 0.9002353549003601
$arr[14,0] = "model_6_4_14"
$arr[14,1] = 0.3113655886044329
$arr[14,2] = -0.01209269523734502
$arr[14,3] = 0.2965313367943251
$arr[14,4] = 0.1985679443675799
$arr[14,5] = 0.7621154189109802
$arr[14,6] = 0.8280742764472961
$arr[14,7] = 0.9821609854698181
$arr[14,8] = 0.9005855917930603
$arr[15,0] = "model_6_4_8"
$arr[15,1] = 0.311393333626768
$arr[15,2] = -0.01285226751600765
$arr[15,3] = 0.2976525559204881
$arr[15,4] = 0.1989307187830328
$arr[15,5] = 0.7620846629142761
$arr[15,6] = 0.8286957740783691
$arr[15,7] = 0.980595588684082
$arr[15,8] = 0.9001778364181519
$arr[16,0] = "model_6_4_16"
$arr[16,1] = 0.3115295285964323
$arr[16,2] = -0.01255186109327644
$arr[16,3] = 0.2976640389539476
$arr[16,4] = 0.1990531777219741
$arr[16,5] = 0.7619339823722839
$arr[16,6] = 0.8284499645233154
$arr[16,7] = 0.9805794954299927
$arr[16,8] = 0.9000402092933655
$arr[17,0] = "model_6_4_13"
$arr[17,1] = 0.3115299047625569
$arr[17,2] = -0.01197225477170627
$arr[17,3] = 0.297410522698823
$arr[17,4] = 0.1991282299230871
$arr[17,5] = 0.7619335055351257
$arr[17,6] = 0.8279757499694824
$arr[17,7] = 0.9809334874153137
$arr[17,8] = 0.8999558687210083
$arr[18,0] = "model_6_4_9"
$arr[18,1] = 0.3118087084660879
$arr[18,2] = -0.01223937936065966
$arr[18,3] = 0.2993575411211055
$arr[18,4] = 0.2001636255902463
$arr[18,5] = 0.7616249322891235
$arr[18,6] = 0.8281943798065186
$arr[18,7] = 0.9782150983810425
$arr[18,8] = 0.8987924456596375
$arr[19,0] = "model_6_4_10"
$arr[19,1] = 0.3118464018691912
$arr[19,2] = -0.01218878596233242
$arr[19,3] = 0.2995263965616289
$arr[19,4] = 0.2002816307850264
$arr[19,5] = 0.7615832090377808
$arr[19,6] = 0.8281528949737549
$arr[19,7] = 0.977979302406311
$arr[19,8] = 0.8986598253250122
$arr[20,0] = "model_6_4_4"
$arr[20,1] = 0.316660198275138
$arr[20,2] = -0.007101853319291251
$arr[20,3] = 0.3259013002828063
$arr[20,4] = 0.2176635751390651
$arr[20,5] = 0.7562557458877563
$arr[20,6] = 0.8239909410476685
$arr[20,7] = 0.9411555528640747
$arr[20,8] = 0.8791274428367615
$arr[21,0] = "model_6_4_3"
$arr[21,1] = 0.3178083047788856
$arr[21,2] = -0.006627919500074242
$arr[21,3] = 0.3325494547054396
$arr[21,4] = 0.2217334403173644
$arr[21,5] = 0.7549852132797241
$arr[21,6] = 0.823603093624115
$arr[21,7] = 0.9318735599517822
$arr[21,8] = 0.8745540380477905
$arr[22,0] = "model_6_4_0"
$arr[22,1] = 0.3224789378909186
$arr[22,2] = 0.016080350249562
$arr[22,3] = 0.4079203220951629
$arr[22,4] = 0.2745547344228872
$arr[22,5] = 0.7498162388801575
$arr[22,6] = 0.8050236701965332
$arr[22,7] = 0.8266431093215942
$arr[22,8] = 0.8151975870132446
$arr[23,0] = "model_6_4_1"
$arr[23,1] = 0.3226135883629595
$arr[23,2] = 0.0129364031071606
$arr[23,3] = 0.3967518755452191
$arr[23,4] = 0.2668131219927011
$arr[23,5] = 0.7496671676635742
$arr[23,6] = 0.80759596824646
$arr[23,7] = 0.8422361612319946
$arr[23,8] = 0.8238970041275024
$arr[24,0] = "model_6_4_2"
$arr[24,1] = 0.3260055067384914
$arr[24,2] = -0.006148803937411795
$arr[24,3] = 0.3923121037541212
$arr[24,4] = 0.2568605495983908
$arr[24,5] = 0.7459133267402649
$arr[24,6] = 0.8232110738754272
$arr[24,7] = 0.8484348654747009
$arr[24,8] = 0.8350809812545776

$ws.Range("A2:I26").Value = $arr
